$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resort rows 48-50: Banglades overtakes Malasia and Panama (col B, total cases) ---
# Row 48 becomes Banglades with its updated totals
$ws.Range("A48").Value = "Banglades"
$ws.Range("B48").Value = 5913
$ws.Range("C48").Value = 497
$ws.Range("D48").Value = 131
$ws.Range("E48").Value = 5630
$ws.Range("F48").Value = 1
$ws.Range("G48").Value = 7
$ws.Range("H48").Value = 152

# Row 49 becomes Malasia with its updated totals
$ws.Range("A49").Value = "Malasia"
$ws.Range("B49").Value = 5820
$ws.Range("C49").Value = 40
$ws.Range("D49").Value = 3957
$ws.Range("E49").Value = 1764
$ws.Range("F49").Value = 37
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 99

# Row 50 becomes Panama, keeping its previous (unchanged) totals
$ws.Range("A50").Value = "Panama"
$ws.Range("B50").Value = 5779
$ws.Range("C50").Value = 241
$ws.Range("D50").Value = 369
$ws.Range("E50").Value = 5245
$ws.Range("F50").Value = 85
$ws.Range("G50").Value = 6
$ws.Range("H50").Value = 165

# --- Row 44 (Noruega): new-case/death tallies updated ---
$ws.Range("E44").Value = 7293
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 202

# --- Row 68 (Uzbekistan): recovered/active counts updated ---
$ws.Range("D68").Value = 804
$ws.Range("E68").Value = 1075

# --- Row 106 (Sri Lanka): recovered/active counts updated ---
$ws.Range("D106").Value = 126
$ws.Range("E106").Value = 390
